$d = $word.ActiveDocument

# 1. Remove the duplicate empty first paragraph at the very start of the
#    document body (an identical empty centered paragraph remains right
#    after it, so the visible content is unchanged apart from one blank
#    line going away).
$first = $d.Paragraphs.Item(1)
$first.Range.Delete()

# 2. Update the header title text.
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("PREWORK de Problemas de regresión ", $true, $false, $false, $false, $false, $true, 1, $false, "Prework - Problemas de regresión ", 2)
